$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.250.46"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.204.36"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'608.06"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").Value = "'156.40"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.203.73"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").Value = "'5.66"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("E12").Value = "  -3.09%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'38.43"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "3.730.32"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "66.374.34"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").Value = "3.202.04"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").Value = "'507.07"
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("D21").Value = "'15.34"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").Value = "'14.63"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("D25").Value = "'85.21"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'3.00"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Value = "'9.04"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("E30").Value = "  +43.73%  "
$ws.Range("D31").Value = "'6.98"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "'28.24"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  -4.98%  "
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'501.35"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'55.41"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").Value = "0.0₃0769"
$ws.Range("E39").Value = "  +13.91%  "
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").Value = "'3.05"
$ws.Range("E42").Value = "  +5.01%  "
$ws.Range("D43").Value = "'8.72"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").Value = "2.910.51"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").Value = "'28.21"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").Value = "'122.22"
$ws.Range("E51").Value = "  +0.24%  "
